# The document contained a long run of empty paragraphs (all sharing the
# same "autoSpaceDE/autoSpaceDN/adjustRightInd" paragraph formatting)
# right after the "- SO 01 Kabelové vedení NN" heading paragraph.
# This edit removes the first 17 of those empty paragraphs, leaving the
# remaining empty paragraphs (and everything else in the document) intact.

$d = $word.ActiveDocument

# Locate the paragraph that contains the heading text.
$findRange = $d.Content
$found = $findRange.Find.Execute("- SO 01 Kabelové vedení NN", $true, $false,
                                  $false, $false, $false, $true, 1, $false,
                                  "", 0)
if (-not $found) {
    Write-Host "ERROR: heading paragraph not found"
}

$headingPara = $findRange.Paragraphs(1)

# The character right after the heading's paragraph mark starts the run
# of empty paragraphs that needs to be trimmed.
$startDel = $headingPara.Range.End

# Walk forward 17 paragraphs (the ones to delete), verifying along the
# way that each is indeed an empty paragraph (just a paragraph mark) so
# we never accidentally remove real content. After the loop, $cur refers
# to the 17th empty paragraph, and its Range.End is the position right
# before the first paragraph that must be kept.
$cur = $headingPara
for ($i = 0; $i -lt 17; $i++) {
    $cur = $cur.Next()
    if ($cur.Range.Text -ne "`r") {
        Write-Host "WARNING: paragraph" $i "after heading is not empty; aborting deletion"
        $startDel = $null
        break
    }
}

if ($startDel -ne $null) {
    $endDel = $cur.Range.End
    $delRange = $d.Range($startDel, $endDel)
    $delRange.Delete()
    Write-Host "Removed 17 empty paragraphs after the SO 01 heading."
}
